$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.217250466346741
$ws.Range("B1").Value = 1.224969863891602
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.007501125335693
$ws.Range("E1").Value = 0.9544448852539062
